$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data to the Gantt chart table
$ws.Range("F9").Value = "V 3"

# Copy the date formatting from the row above so the new date cell
# matches the existing short-date style used in the table
$ws.Range("G8").Copy()
$ws.Range("G9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("G9").Value = 43997

$ws.Range("H9").Value = 5
$ws.Range("I9").Value = "Fixed Dates And Added Unit 4"

# Update selection to match the authored state
$ws.Range("P7").Select()
